$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update price (D) and volume-change (E) columns, and for the two rank swaps
# also update coin name (B) and link (C) columns, to match the refreshed data pull.

$ws.Range("D2").Value = "59.096.12"
$ws.Range("E2").Value = "  +1.78%  "
$ws.Range("D3").Value = "2.543.18"
$ws.Range("E3").Value = "  +3.03%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "527.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.54%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.29"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.49%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.38%  "
$ws.Range("E8").Value = "  +1.63%  "
$ws.Range("D9").Value = "2.540.89"
$ws.Range("E9").Value = "  +2.49%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0991"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.96%  "
$ws.Range("E11").Value = "  -1.38%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.19"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.90%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.336"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.13%  "
$ws.Range("D14").Value = "2.992.20"
$ws.Range("E14").Value = "  +2.86%  "
$ws.Range("D15").Value = "59.020.89"
$ws.Range("E15").Value = "  +1.75%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "22.41"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.19%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000136"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.26%  "
$ws.Range("D18").Value = "2.541.88"
$ws.Range("E18").Value = "  +2.98%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.76"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.61%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "324.45"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.66%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.21"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.26%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.13"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +7.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.34"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.15%  "
$ws.Range("E25").Value = "  +0.46%  "
$ws.Range("B26").Value = "Binance-PegBSC-USD"
$ws.Range("C26").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.997"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.06%  "
$ws.Range("B27").Value = "Kaspa"
$ws.Range("C27").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.161"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.13%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.45"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.73%  "
$ws.Range("D29").Value = "0.0₃0758"
$ws.Range("E29").Value = "  +1.71%  "
$ws.Range("E30").Value = "  +3.40%  "
$ws.Range("E31").Value = "  +2.63%  "
$ws.Range("B32").Value = "Monero"
$ws.Range("C32").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "168.70"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.58%  "
$ws.Range("B33").Value = "Aptos"
$ws.Range("C33").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.38"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.82%  "
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.24%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.35"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.87%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.27"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.65%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.98"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.03%  "
$ws.Range("E39").Value = "  +2.76%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.72"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.31%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.788"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.77%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "282.52"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.33%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.49"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.85%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.08"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.29%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "131.85"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +6.79%  "
$ws.Range("E46").Value = "  +2.17%  "
$ws.Range("E47").Value = "  +1.63%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0507"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.97%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "17.88"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.73%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0217"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.04%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.22"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.58%  "
